$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 658, shifting existing rows 658:720 down to 659:721
$ws.Rows("658:658").Insert()

# Populate the newly inserted row 658 with its data
$ws.Range("A658").Value = 4
$ws.Range("B658").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C658").Value = "Los Lagos"
$ws.Range("D658").Value = 45132
$ws.Range("E658").Value = 10
$ws.Range("F658").Value = 100114001
$ws.Range("G658").Value = "Papa"
$ws.Range("H658").Value = "Patagonia"
$ws.Range("I658").Value = "1a (guarda)"
$ws.Range("J658").Value = 500
$ws.Range("K658").Value = 18000
$ws.Range("L658").Value = 18000
$ws.Range("M658").Value = 18000
$ws.Range("N658").Value = "$/saco 25 kilos"
$ws.Range("O658").Value = "Provincia de Llanquihue"
$ws.Range("P658").Value = 720
$ws.Range("Q658").Value = 25
$ws.Range("R658").Value = "Hortaliza"
